# Menu.xlsx - "fixed: outputting specials and exclusive"
# Adds two new dish rows (blueberry, ice cream) to the food menu sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=DishName, B=Price, C=ChefSpecial, D=ValentineExclusive

# Row 13: blueberry / 56 / not a chef special / not a valentine exclusive
$ws.Cells.Item(13, 1).Value = "blueberry"
$ws.Cells.Item(13, 2).Value = 56
$ws.Cells.Item(13, 3).Value = $false
$ws.Cells.Item(13, 4).Value = $false

# Row 14: ice cream / 44 / chef special / valentine exclusive
$ws.Cells.Item(14, 1).Value = "ice cream"
$ws.Cells.Item(14, 2).Value = 44
$ws.Cells.Item(14, 3).Value = $true
$ws.Cells.Item(14, 4).Value = $true
